# Registro de atividades.xlsx - add row 10 entry (4th revision) to the
# activity log table on Plan1 / "Registro de atualizações".
#
# Target state (per the authoritative diff):
#   - B10 = 30/11/2018 (serial 43434), short-date format
#   - C10 = "Resumo(revisado), Projeto final, Plano de testes,
#            Dificuldades encontradas, manual do usuario(revisado)"
#            (new shared string), wrapped text
#   - D10 = "Lucas Armando Ciello" (same author as the other rows)
#   - row 10 grows taller to fit the wrapped summary text
#   - selection cursor left on J9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the new revision row -----------------------------------------

# Date of the revision (2018-11-30). Assign the numeric serial directly so
# the cell is stored as a real date number, then apply the same short-date
# display the other date cells in the column use (builtin format id 14,
# "m/d/yyyy" once localized - expressed here with its canonical ECMA code so
# it maps onto the existing builtin format instead of a custom one).
$ws.Range("B10").Value = 43434
$ws.Range("B10").NumberFormat = "mm-dd-yy"

# Summary text for this revision; wrap so the whole cell content is visible.
$ws.Range("C10").Value = "Resumo(revisado), Projeto final, Plano de testes, Dificuldades encontradas, manual do usuario(revisado)"
$ws.Range("C10").WrapText = $true

# Same author as the previous rows.
$ws.Range("D10").Value = "Lucas Armando Ciello"

# Grow row 10 so the wrapped summary text is fully visible.
$ws.Rows.Item(10).RowHeight = 30.75

# --- Leave the selection where the editor left it --------------------------
$ws.Range("J9").Select()
